$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Price (column D) updates ---
Set-TextValue "D2" "42.839.94"
Set-TextValue "D3" "2.210.19"
Set-TextValue "D5" "256.62"
Set-TextValue "D7" "77.38"
Set-TextValue "D9" "0.593"
Set-TextValue "D10" "42.81"
Set-TextValue "D11" "0.0910"
Set-TextValue "D12" "6.98"
Set-TextValue "D13" "0.102"
Set-TextValue "D14" "2.543.61"
Set-TextValue "D15" "14.41"
Set-TextValue "D16" "2.215.92"
Set-TextValue "D17" "0.782"
Set-TextValue "D18" "42.800.87"
Set-TextValue "D20" "71.01"
Set-TextValue "D21" "5.97"
Set-TextValue "D23" "229.77"
Set-TextValue "D24" "9.20"
Set-TextValue "D26" "42.57"
Set-TextValue "D27" "10.71"
Set-TextValue "D30" "2.21"
Set-TextValue "D31" "173.30"
Set-TextValue "D32" "20.35"
Set-TextValue "D33" "0.0875"
Set-TextValue "D34" "5.20"
Set-TextValue "D36" "0.0356"
Set-TextValue "D38" "4.38"
Set-TextValue "D39" "13.05"
Set-TextValue "D41" "2.10"
Set-TextValue "D42" "0.201"
Set-TextValue "D43" "60.74"
Set-TextValue "D44" "5.30"
Set-TextValue "D45" "102.93"
Set-TextValue "D48" "0.0970"
Set-TextValue "D49" "1.11"

# --- Volume(1h) (column E) updates ---
Set-TextValue "E2" "  -0.54%  "
Set-TextValue "E3" "  -1.31%  "
Set-TextValue "E5" "  +2.53%  "
Set-TextValue "E6" "  +0.51%  "
Set-TextValue "E7" "  +3.11%  "
Set-TextValue "E8" "  -0.02%  "
Set-TextValue "E9" "  -0.94%  "
Set-TextValue "E10" "  +4.53%  "
Set-TextValue "E11" "  -2.20%  "
Set-TextValue "E12" "  +1.19%  "
Set-TextValue "E13" "  +0.71%  "
Set-TextValue "E14" "  -1.25%  "
Set-TextValue "E15" "  -1.51%  "
Set-TextValue "E16" "  -2.58%  "
Set-TextValue "E17" "  -1.15%  "
Set-TextValue "E18" "  -0.42%  "
Set-TextValue "E19" "  -1.13%  "
Set-TextValue "E20" "  -0.17%  "
Set-TextValue "E21" "  -0.31%  "
Set-TextValue "E22" "  +4.68%  "
Set-TextValue "E23" "  +0.17%  "
Set-TextValue "E24" "  -5.08%  "
Set-TextValue "E26" "  +7.70%  "
Set-TextValue "E27" "  -0.58%  "
Set-TextValue "E28" "  -2.59%  "
Set-TextValue "E29" "  -2.53%  "
Set-TextValue "E30" "  -1.53%  "
Set-TextValue "E31" "  +1.06%  "
Set-TextValue "E32" "  +0.69%  "
Set-TextValue "E33" "  +9.51%  "
Set-TextValue "E34" "  -0.95%  "
Set-TextValue "E35" "  +0.01%  "
Set-TextValue "E36" "  +8.10%  "
Set-TextValue "E37" "  -2.88%  "
Set-TextValue "E38" "  -2.21%  "
Set-TextValue "E39" "  +0.61%  "
Set-TextValue "E40" "  +17.58%  "
Set-TextValue "E41" "  -0.37%  "
Set-TextValue "E42" "  -2.13%  "
Set-TextValue "E43" "  +2.19%  "
Set-TextValue "E44" "  -2.40%  "
Set-TextValue "E45" "  -0.73%  "
Set-TextValue "E48" "  -2.09%  "
Set-TextValue "E49" "  +0.79%  "
Set-TextValue "E50" "  -1.40%  "
Set-TextValue "E51" "  +22.75%  "

# --- Rows 46 & 47 swap: WOONetwork/FraxShare order reversed, with updated price/volume ---
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D46" "8.45"
Set-TextValue "E46" "  -2.24%  "

$ws.Range("B47").Value = "WOONetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextValue "D47" "0.472"
Set-TextValue "E47" "  -2.19%  "
